$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Level" column (C) values:
# Row 3 (grandmother): moderate -> easy
# Row 5 (grandfather): easy -> hard
# Row 6 (brother): easy -> hard
$ws.Range("C3").Value = "easy"
$ws.Range("C5").Value = "hard"
$ws.Range("C6").Value = "hard"
